# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the author's fix: the scraper previously pulled only team
# statistics, not the season win/loss/tie record, so these three columns
# are appended after the existing "Unnamed: 28" column (AC) and populated
# for every player row with this team's season record: 51-111-0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cells the same formatting as the existing header
# row (bold, centered, bordered) by copying the style of the last header
# cell (AC1) onto the new ones before writing their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-52) gets the same team season record.
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 51   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 111  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
